$d = $word.ActiveDocument

# 1) Objective paragraph: "...such as EVO to help..." -> "...such as EVO and Trimble to help..."
$d.Content.Find.Execute("EVO to help create", $true, $false, $false, $false, $false, $true, 1, $false, "EVO and Trimble to help create", 2) | Out-Null

# 2) First job entry date range: "2020 - Current " -> "2019 - 2022 "
$d.Content.Find.Execute("2020 - Current ", $true, $false, $false, $false, $false, $true, 1, $false, "2019 - 2022 ", 2) | Out-Null

# 3) First job entry title: "Web Developer" -> "Front End Developer"
$d.Content.Find.Execute("Web Developer", $true, $false, $false, $false, $false, $true, 1, $false, "Front End Developer", 2) | Out-Null

# 4) First job entry company: "Self Employed" -> "Trimble Transportation"
$d.Content.Find.Execute("Self Employed", $true, $false, $false, $false, $false, $true, 1, $false, "Trimble Transportation", 2) | Out-Null

# 5) First job entry description
$oldDesc = "Designed and created web applications for clients utilizing frameworks such as React for reusable and scalable code or Shopify to create an online store. Worked closely with clients to implement their ideas to create and manage web apps."
$newDesc = "Designed and implemented web applications for clients using Trimble app by creating reusable and scalable code. Worked closely with team to implement ideas to create and manage web app."
$d.Content.Find.Execute($oldDesc, $true, $false, $false, $false, $false, $true, 1, $false, $newDesc, 2) | Out-Null

# 6) Second job entry (EVO Payments) date range: "2019 - 2020 " -> "2017 - 2019 "
$d.Content.Find.Execute("2019 - 2020 ", $true, $false, $false, $false, $false, $true, 1, $false, "2017 - 2019 ", 2) | Out-Null

# 7) Education date range: "2012-2014" -> "2012-2016"
$d.Content.Find.Execute("2012-2014", $true, $false, $false, $false, $false, $true, 1, $false, "2012-2016", 2) | Out-Null

# 8) Education degree: "Associates" -> "Bachelors"
$d.Content.Find.Execute("Associates in Computer Science", $true, $false, $false, $false, $false, $true, 1, $false, "Bachelors in Computer Science", 2) | Out-Null

# 9) Education institution: "Front Range Community College" -> "Colorado State University"
$d.Content.Find.Execute("Front Range Community College", $true, $false, $false, $false, $false, $true, 1, $false, "Colorado State University", 2) | Out-Null
